$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.486.03"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "1.571.44"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.88"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.10"
$ws.Range("E8").Value = "  +5.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.10"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.248"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0882"
$ws.Range("D13").Value = "1.795.23"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "1.571.19"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.69"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "28.454.97"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.11"
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.18"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.89"
$ws.Range("E23").Value = "  -5.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("E25").Value = "  +5.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.06"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.98"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "1.390.39"
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("E39").Value = "  +3.65%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.59"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.98"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "1.708.14"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.86"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("E51").Value = "  -1.30%  "
